$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 93, shifting the old (empty) row 93 and the
# summary-formula row 94 down by one.
$ws.Rows("93").Insert()

# Fill in the new journal entry.
$ws.Range("A93").Value = 44209
$ws.Range("B93").Value = "Modification des pages contact, intelligence_artificielle et robotique"
$ws.Range("C93").Value = "Insertion des mentions légales des logos en footer"
$ws.Range("D93").Value = "Mathieu"
$ws.Range("E93").Value = "VS Code"
$ws.Range("F93").Value = "HTML"
$ws.Rows("93").RowHeight = 30

# Update the summary formula (now on row 95) to include the new row.
$ws.Range("D95").Formula = "=COUNTIF(D1:D93,""Mathieu"")/92*100"

# Extend the AutoFilter range to cover the new data.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:F93").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the AutoFilter range.
$n = $wb.Names.Item(1)
$n.RefersTo = "=Feuil1!`$A`$1:`$F`$93"

# Restore the selection to where it ends up in the saved workbook.
[void]$ws.Range("D96").Select()
